# Weekly update: insert a new price record as row 40, pushing the
# existing rows 40-162 down to 41-163 (dimension grows from A1:R162 to
# A1:R163).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40..162 down to 41..163, leaving a blank row 40 to fill in.
$ws.Rows(40).Insert()

# Populate the newly inserted row 40 with this week's record.
$ws.Range("A40").Value = 8
$ws.Range("B40").Value = "Terminal La Palmera de La Serena"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44608
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = 100112037
$ws.Range("G40").Value = "Cebollín"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 600
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = 1100
$ws.Range("N40").Value = "$/paquete 6 unidades"
$ws.Range("O40").Value = "Provincia del Elquí"
$ws.Range("P40").Value = 183
$ws.Range("Q40").Value = 6
$ws.Range("R40").Value = "Hortaliza"
